# Weekly update to the Ají (Vega Central Mapocho de Santiago) price sheet.
# A new weekly price observation is inserted as row 236 (pushing the
# existing rows 236:247 down to 237:248), matching the new used range
# A1:R248 (was A1:R247).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 236:247 down by inserting a fresh row at 236. Excel carries
# the formatting of the row above down onto the new row (same behaviour
# as inserting a row in the desktop app), which reproduces the D-column
# date style (s="2") on the new row automatically.
$ws.Rows(236).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(236, 1).Value2  = 9
$ws.Cells.Item(236, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(236, 3).Value2  = "Metropolitana"
$ws.Cells.Item(236, 4).Value2  = 44610
$ws.Cells.Item(236, 5).Value2  = 13
$ws.Cells.Item(236, 6).Value2  = 100112021
$ws.Cells.Item(236, 7).Value2  = "Ají"
$ws.Cells.Item(236, 8).Value2  = "Americana (o)"
$ws.Cells.Item(236, 9).Value2  = "Primera"
$ws.Cells.Item(236, 10).Value2 = 26
$ws.Cells.Item(236, 11).Value2 = 16000
$ws.Cells.Item(236, 12).Value2 = 18000
$ws.Cells.Item(236, 13).Value2 = 17000
$ws.Cells.Item(236, 14).Value2 = "`$/caja 25 kilos"
$ws.Cells.Item(236, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(236, 16).Value2 = 680
$ws.Cells.Item(236, 17).Value2 = 25
$ws.Cells.Item(236, 18).Value2 = "Hortaliza"
